$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 18:22"

# --- Country ranking refresh: Argelia overtakes Egipto ---
$ws.Range("A57").Value = "Argelia"
$ws.Range("B57").Value = 1825
$ws.Range("C57").Value = 64
$ws.Range("D57").Value = 460
$ws.Range("E57").Value = 1090
$ws.Range("F57").Value = 46
$ws.Range("G57").Value = 19
$ws.Range("H57").Value = 275
$ws.Range("A58").Value = "Egipto"
$ws.Range("B58").Value = 1794
$ws.Range("D58").Value = 384
$ws.Range("E58").Value = 1275
$ws.Range("F58").Value = 0
$ws.Range("H58").Value = 135

# --- Country ranking refresh: Irak overtakes Nueva Zelanda, Hungria, Estonia ---
$ws.Range("A63").Value = "Irak"
$ws.Range("B63").Value = 1318
$ws.Range("C63").Value = 39
$ws.Range("D63").Value = 601
$ws.Range("E63").Value = 645
$ws.Range("F63").Value = 0
$ws.Range("H63").Value = 72
$ws.Range("A64").Value = "Nueva Zelanda"
$ws.Range("B64").Value = 1312
$ws.Range("C64").Value = 29
$ws.Range("D64").Value = 422
$ws.Range("E64").Value = 886
$ws.Range("F64").Value = 5
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 4
$ws.Range("A65").Value = "Hungria"
$ws.Range("B65").Value = 1310
$ws.Range("C65").Value = 120
$ws.Range("D65").Value = 115
$ws.Range("E65").Value = 1110
$ws.Range("F65").Value = 17
$ws.Range("G65").Value = 8
$ws.Range("H65").Value = 85
$ws.Range("A66").Value = "Estonia"
$ws.Range("B66").Value = 1304
$ws.Range("C66").Value = 46
$ws.Range("D66").Value = 93
$ws.Range("E66").Value = 1187
$ws.Range("F66").Value = 11
$ws.Range("H66").Value = 24

# --- Remaining numeric refreshes (country order unchanged) ---
$ws.Range("B4").Value = 507338
$ws.Range("C4").Value = 4462
$ws.Range("D4").Value = 28224
$ws.Range("E4").Value = 459399
$ws.Range("G4").Value = 968
$ws.Range("H4").Value = 19715
$ws.Range("B6").Value = 152271
$ws.Range("C6").Value = 4694
$ws.Range("D6").Value = 32534
$ws.Range("E6").Value = 100269
$ws.Range("F6").Value = 3381
$ws.Range("G6").Value = 619
$ws.Range("H6").Value = 19468
$ws.Range("B19").Value = 13795
$ws.Range("C19").Value = 235
$ws.Range("E19").Value = 6854
$ws.Range("D35").Value = 411
$ws.Range("E35").Value = 5291
$ws.Range("G35").Value = 10
$ws.Range("H35").Value = 129
$ws.Range("D85").Value = 77
$ws.Range("E85").Value = 522
$ws.Range("F85").Value = 34
$ws.Range("B110").Value = 263
$ws.Range("C110").Value = 6
$ws.Range("E110").Value = 256
$ws.Range("B112").Value = 242
$ws.Range("C112").Value = 8
$ws.Range("E112").Value = 183
$ws.Range("B162").Value = 31
$ws.Range("C162").Value = 4
$ws.Range("E162").Value = 26
